# Update p-values (column F, "Pr(>F)") on Sheet 1 to values rounded to
# three/four digits as produced by the re-run PERMANOVA tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$updates = @{
    "F2"  = 0.4525
    "F3"  = 0.0113
    "F4"  = 0.4617
    "F7"  = 0.1631
    "F8"  = 0.0494
    "F9"  = 0.9533
    "F12" = 0.3356
    "F13" = 0.5926
    "F14" = 0.4014
    "F17" = 0.0737
    "F18" = 0.9975
    "F19" = 0.5746
    "F22" = 0.183
    "F23" = 0.4906
    "F24" = 0.2146
    "F27" = 0.8906
    "F28" = 0.0138
    "F29" = 0.8931
    "F32" = 0.1669
    "F33" = 0.7037
    "F34" = 0.9623
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
